$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the existing header row (row 1 -> row 3),
# shifting the header and all data rows down by two.
$ws.Rows("1:2").Insert() | Out-Null

# The two hidden "extra" rows (formerly rows 12:13, now rows 14:15 after the
# insert above) are removed entirely rather than merely unhidden.
$ws.Rows("14:15").Delete() | Out-Null

# Set up best-fit-like column widths for the data columns (A:E).
$ws.Columns(1).ColumnWidth = 3.1666666666666665
$ws.Columns(2).ColumnWidth = 9.666666666666666
$ws.Columns(3).ColumnWidth = 9.333333333333334
$ws.Columns(4).ColumnWidth = 9
$ws.Columns(5).ColumnWidth = 9

# Move the active selection.
$ws.Range("I7").Select() | Out-Null
